# Apply the latest crypto price/volume snapshot to the 'cryptos' worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price column values look like plain decimal numbers (e.g. 533.38).
# Prefixing with an apostrophe forces Excel to store them as text (quote-prefixed),
# matching the original workbook where every Price cell is a text string.
$apos = "'"

$ws.Range("D2").Value = '58.374.91'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '3.141.78'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = $apos + '533.38'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = $apos + '142.63'
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.140.41'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").Value = $apos + '0.446'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("D12").Value = $apos + '0.393'
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("D13").Value = '3.682.77'
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("D15").Value = $apos + '25.66'
$ws.Range("E15").Value = '  -4.72%  '
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '58.404.90'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '3.138.96'
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").Value = $apos + '12.82'
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("D22").Value = $apos + '343.85'
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = $apos + '67.70'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = $apos + '0.995'
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("D28").Value = '0.0₃0936'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = $apos + '7.43'
$ws.Range("E29").Value = '  +2.41%  '
$ws.Range("D31").Value = $apos + '6.42'
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("E32").Value = '  +1.47%  '
$ws.Range("D33").Value = $apos + '21.12'
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("E35").Value = '  +3.11%  '
$ws.Range("D36").Value = $apos + '158.05'
$ws.Range("E36").Value = '  +2.40%  '
$ws.Range("D37").Value = $apos + '6.24'
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("D38").Value = $apos + '26.34'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("E39").Value = '  -4.76%  '
$ws.Range("D40").Value = $apos + '1.67'
$ws.Range("E40").Value = '  +12.08%  '
$ws.Range("D41").Value = $apos + '0.0668'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("D42").Value = $apos + '0.708'
$ws.Range("E42").Value = '  +5.06%  '
$ws.Range("D43").Value = $apos + '4.01'
$ws.Range("E43").Value = '  +2.55%  '
$ws.Range("D44").Value = '3.182.86'
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("D45").Value = $apos + '36.71'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("E47").Value = '  +2.73%  '
$ws.Range("D48").Value = '2.279.82'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("E49").Value = '  +4.57%  '
$ws.Range("D50").Value = $apos + '20.70'
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("E51").Value = '  +1.51%  '

Write-Output "done"
